$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7061076164245605
$ws.Range("B1").Value = 1.424153685569763
$ws.Range("C1").Value = 4.145806789398193
$ws.Range("D1").Value = 2.486992835998535
$ws.Range("E1").Value = 0.5606345534324646
